$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new week of data (2022-03-17, serial 44637) is inserted at
# the top of the Apio time series (row 206), pushing the existing rows 206-251
# down to 208-253.
$ws.Rows.Item(206).EntireRow.Insert()
$ws.Rows.Item(206).EntireRow.Insert()

# New row 206 - "Primera" quality, new week
$ws.Cells.Item(206, 1).Value = 9
$ws.Cells.Item(206, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(206, 3).Value = "Metropolitana"
$ws.Cells.Item(206, 4).Value = 44637
$ws.Cells.Item(206, 5).Value = 13
$ws.Cells.Item(206, 6).Value = 100112017
$ws.Cells.Item(206, 7).Value = "Apio"
$ws.Cells.Item(206, 8).Value = "Americana (o)"
$ws.Cells.Item(206, 9).Value = "Primera"
$ws.Cells.Item(206, 10).Value = 79
$ws.Cells.Item(206, 11).Value = 9000
$ws.Cells.Item(206, 12).Value = 10000
$ws.Cells.Item(206, 13).Value = 9494
$ws.Cells.Item(206, 14).Value = '$/docena de matas'
$ws.Cells.Item(206, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(206, 16).Value = 1582
$ws.Cells.Item(206, 17).Value = 6
$ws.Cells.Item(206, 18).Value = "Hortaliza"

# New row 207 - "Segunda" quality, new week
$ws.Cells.Item(207, 1).Value = 9
$ws.Cells.Item(207, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(207, 3).Value = "Metropolitana"
$ws.Cells.Item(207, 4).Value = 44637
$ws.Cells.Item(207, 5).Value = 13
$ws.Cells.Item(207, 6).Value = 100112017
$ws.Cells.Item(207, 7).Value = "Apio"
$ws.Cells.Item(207, 8).Value = "Americana (o)"
$ws.Cells.Item(207, 9).Value = "Segunda"
$ws.Cells.Item(207, 10).Value = 43
$ws.Cells.Item(207, 11).Value = 8000
$ws.Cells.Item(207, 12).Value = 8000
$ws.Cells.Item(207, 13).Value = 8000
$ws.Cells.Item(207, 14).Value = '$/docena de matas'
$ws.Cells.Item(207, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(207, 16).Value = 1333
$ws.Cells.Item(207, 17).Value = 6
$ws.Cells.Item(207, 18).Value = "Hortaliza"
